$wb = $excel.ActiveWorkbook

# --- Update "Last Updated" timestamp on the Metadata sheet ---
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("A2").Value = "05 Nov 2025, 12:50 PM"

# --- Update the "1 Year" return column (F2:F76) on the Industry Analysis sheet ---
$ws = $wb.Worksheets.Item("Industry Analysis")

$updates = @{
    "F2" = 21.0016
    "F3" = -16.2396
    "F4" = 27.1317
    "F5" = -50.6494
    "F6" = 53.2813
    "F7" = -8.106199999999999
    "F8" = -9.552099999999999
    "F9" = 36.3756
    "F10" = -6.1314
    "F11" = 31.9081
    "F12" = -18.4955
    "F13" = 14.0155
    "F14" = -36.0718
    "F15" = -0.1622
    "F16" = 0.1459
    "F17" = -22.0012
    "F18" = 1.0561
    "F19" = -27.708
    "F20" = 47.7309
    "F21" = 12.0959
    "F22" = 95.1491
    "F23" = -50.2657
    "F24" = -13.3427
    "F25" = -9.9316
    "F26" = 5.8244
    "F27" = -32.7692
    "F28" = -24.8224
    "F29" = -18.4191
    "F30" = 25.8569
    "F31" = 58.4712
    "F32" = -3.3862
    "F33" = -6.3282
    "F34" = 27.7203
    "F35" = 4.4873
    "F36" = -4.9458
    "F37" = 3.6074
    "F38" = -23.3973
    "F39" = 8.7355
    "F40" = -5.8541
    "F41" = -8.3934
    "F42" = 20.3818
    "F43" = 14.3164
    "F44" = -12.6846
    "F45" = 28.4075
    "F46" = -1.1135
    "F47" = -37.1997
    "F48" = -29.8569
    "F49" = -27.5511
    "F50" = -49.7478
    "F51" = -51.8002
    "F52" = -38.5254
    "F53" = -12.4886
    "F54" = -5.0725
    "F55" = -17.7445
    "F56" = -26.636
    "F57" = -29.3361
    "F58" = -11.9574
    "F59" = -24.5687
    "F60" = -12.3
    "F61" = -10.9446
    "F62" = -17.1229
    "F63" = -9.5038
    "F64" = 54.2749
    "F65" = -43.4736
    "F66" = 13.2687
    "F67" = 12.7149
    "F68" = 24.8057
    "F69" = -17.0328
    "F70" = -6.8927
    "F71" = 13.6034
    "F72" = 3.9995
    "F73" = -16.226
    "F74" = -16.2448
    "F75" = 28.6924
    "F76" = 48.9752
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
